# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.088.58"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.143.63"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.10"
$ws.Range("E5").Value = "  +9.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "636.42"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("E7").Value = "  +2.79%  "
$ws.Range("E8").Value = "  -3.13%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.141.32"
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.723"
$ws.Range("E11").Value = "  -5.36%  "
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.97"
$ws.Range("E13").Value = "  +5.48%  "
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.947.21"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.721.48"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.169.35"
$ws.Range("E18").Value = "  +2.65%  "
$ws.Range("E19").Value = "  -2.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.33"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "446.37"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  +8.04%  "
$ws.Range("E24").Value = "  +2.10%  "
$ws.Range("E25").Value = "  -4.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.19"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.54"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.71"
$ws.Range("E30").Value = "  +5.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.160"
$ws.Range("E31").Value = "  -3.93%  "
$ws.Range("E32").Value = "  +8.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.12"
$ws.Range("E33").Value = "  +11.23%  "
$ws.Range("E34").Value = "  +22.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.83"
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "515.73"
$ws.Range("E36").Value = "  -2.59%  "
$ws.Range("E37").Value = "  +1.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.21"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  +2.96%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.418"
$ws.Range("E41").Value = "  +2.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.21"
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0852"
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.35"
$ws.Range("E45").Value = "  +48.90%  "
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.74"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.696"
$ws.Range("E48").Value = "  +9.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.83"
$ws.Range("E49").Value = "  +3.31%  "
$ws.Range("E50").Value = "  +3.46%  "
$ws.Range("E51").Value = "  +4.84%  "
